$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Wraps a fragment of <w:p>...</w:p> markup (using only the bare "w:" prefix,
# Word resolves it against the real document namespace on insertion) in the
# minimal WordprocessingML package envelope that Range.InsertXML expects.
function New-OpenXmlPackage($bodyFragment) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $bodyFragment + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Inserts a brand-new, literal <w:p>...</w:p> paragraph (exactly as given in
# $paragraphXml, no inherited pStyle/pPr) right after the first paragraph
# whose full text equals $afterText. Splicing happens immediately before
# that paragraph's own trailing paragraph mark, which is the one spot where
# InsertXML adds content as a clean sibling <w:p> rather than folding it
# into neighbouring runs.
function Insert-ParagraphAfter($afterText, $paragraphXml) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd("`r", "`a", "`n")
        if ($t -eq $afterText) {
            $insertPos = $p.Range.End - 1
            $r = $d.Range($insertPos, $insertPos)
            $r.InsertXML((New-OpenXmlPackage ('<w:body>' + $paragraphXml + '</w:body>')))
            return $true
        }
    }
    throw "Insert-ParagraphAfter: paragraph with text [$afterText] not found"
}

# Replaces the body of the first paragraph whose full text equals $matchText
# with the literal run markup in $runsXml (exact <w:r>...</w:r> runs, no
# normalisation/merging of adjacent identically-formatted runs).
function Set-ParagraphRuns($matchText, $runsXml) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text.TrimEnd("`r", "`a", "`n")
        if ($t -eq $matchText) {
            $bodyRange = $d.Range($p.Range.Start, $p.Range.End - 1)
            $bodyRange.InsertXML((New-OpenXmlPackage ('<w:body>' + $runsXml + '</w:body>')))
            return $true
        }
    }
    throw "Set-ParagraphRuns: paragraph with text [$matchText] not found"
}

# ---------------------------------------------------------------------------
# 1. Title: "LOM3250 -  Trabalho de Graduação" -> "... II"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("LOM3250 -  Trabalho de Graduação", $true, $false, $false, $false, $false, `
    $true, 1, $false, "LOM3250 -  Trabalho de Graduação II", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Subtitle: "Undergraduate Work" -> "Graduation Monograph II"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Undergraduate Work", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Graduation Monograph II", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Activation date: 01/01/2012 -> 01/01/2023
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Ativação: 01/01/2012", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Ativação: 01/01/2023", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. New italic English paragraph after the "Objetivos" body paragraph
# ---------------------------------------------------------------------------
$objetivosBody = "O Trabalho de Graduação (TG) tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de tarefas que fazem parte do perfil de atuação profissional do engenheiro físico."
$objetivosEn = '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t>The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student''s ability to perform tasks that are part of the professional performance profile of the physical engineer.</w:t></w:r></w:p>'
Insert-ParagraphAfter $objetivosBody $objetivosEn

# ---------------------------------------------------------------------------
# 5. Professor list: single entry -> two entries (separate runs, line break
#    only after the first one)
# ---------------------------------------------------------------------------
$professorsXml = '<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>5840730 - Antonio Jefferson da Silva Machado</w:t><w:br/></w:r><w:r><w:t>1176388 - Luiz Tadeu Fernandes Eleno</w:t></w:r></w:p>'
Set-ParagraphRuns "519033 - Carlos Yujiro Shigue" $professorsXml

# ---------------------------------------------------------------------------
# 6. New italic English paragraph after the "Programa resumido" body
# ---------------------------------------------------------------------------
$resumidoBody = "Elaborar uma monografia de Trabalho de Graduação sob a orientação de docente e apresentá-la perante uma banca de examinadores."
$resumidoEn = '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t>Prepare a monograph of Undergraduate Work under the guidance of a professor and present it to a panel of examiners.</w:t></w:r></w:p>'
Insert-ParagraphAfter $resumidoBody $resumidoEn

# ---------------------------------------------------------------------------
# 7. New italic English paragraph after the "Programa" body
# ---------------------------------------------------------------------------
$programaBody = "O programa da disciplina será constituído pelas seguintes etapas: 1) Propor no início do período letivo um plano de trabalho a ser avaliado por uma comissão de professores. 2) Elaborar a monografia cujo tema seja pertencente ao conteúdo programático do curso de Engenharia Física, podendo ser um tópico de interesse técnico ou científico, estudo de caso ou uma proposta de projeto. 3) Definição e divulgação da data de apresentação após a entrega da monografia com antecedência de, no mínimo, 15 dias úteis. 4) Definição da banca de examinadores, sendo constituída pelo professor orientador e por no mínimo dois professores convidados. 5) Apresentação e avaliação do TG. 6) Divulgação da avaliação. Em caso de aprovação, deverá ser feita a entrega do exemplar final da monografia (cópia impressa e eletrônica) com o de acordo do professor orientador."
$programaEn = '<w:p><w:r><w:rPr><w:i/></w:rPr><w:t>The course program will consist of the following steps: 1) Preparation and writing of a monograph on a previously defined and approved subject in the Undergraduate Work I discipline. 2) Definition and disclosure of the presentation date after delivery of the monograph in advance of at least , 15 working days. 3) Definition of the panel of examiners, consisting of the supervisor and at least two invited professionals, with training in engineering or related areas. 4) Presentation and evaluation of the TG. 5) Publication of the evaluation. In case of approval, the final copy of the monograph (printed and electronic copy) must be delivered with the agreement of the supervisor.</w:t></w:r></w:p>'
Insert-ParagraphAfter $programaBody $programaEn

# ---------------------------------------------------------------------------
# 8. Requisito: LOM3238 Projeto Integrado I -> LOM3267 Trabalho de Graduação I
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("LOM3238 -  Projeto Integrado I  (Requisito)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "LOM3267 -  Trabalho de Graduação I  (Requisito)", 2) | Out-Null

Write-Output "All edits applied successfully"
